$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column H (8), shifting H:I -> J:K
$ws.Range("H1:I1").EntireColumn.Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# New header cells with same style as other headers (copy from G4)
$ws.Range("G4").Copy()
$ws.Range("H4:I4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("H4").Value = "Kredit Dus"
$ws.Range("I4").Value = "Kredit Pack"

# Column widths for the new H:I columns (match column G's width of 15.42578125
# characters as closely as the ColumnWidth pixel-grid allows)
$ws.Range("H1:I1").EntireColumn.ColumnWidth = 14.67

# Update selection to I4
$ws.Range("I4").Select()
